$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "hey"
$ws.Range("F1").Value = "top"

$ws.Range("F1").Select()
